# Daily_Updates.xlsx - 2025-09-08 OneDrive sync edit
# Inserts a new "Sheet1" worksheet between "Daily Updates" and "Sections Config"
# containing a small status/handover table, then makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new worksheet in the right spot -------------------------
$sectionsConfig = $wb.Worksheets.Item("Sections Config")
$new = $wb.Worksheets.Add($sectionsConfig)
$new.Name = "Sheet1"

# --- 2. Fill in the table (B2:E7) ------------------------------------------
$new.Range("B2").Value = "WORK"
$new.Range("C2").Value = "FINAL DATE"
$new.Range("D2").Value = "LATEST UPDATE"
$new.Range("E2").Value = "SUPPORT"

$new.Range("B3").Value = "VALLAM"
$new.Range("C3").Value = "Approx 30/11/2025"
$new.Range("D3").Value = "G12R Module being sent from Falta plant and M10R being sent from Chennai Plant" + [char]10 + "Modules will be dispatch by 10/09/2025 (Pending SO in Falta custom)"

$new.Range("B4").Value = "DEVELOPMENT OF QC"
$new.Range("C4").Value = 45945
$new.Range("D4").Value = "MAILED TO QC FOR SENDING THE SAMPLES"
$new.Range("E4").Value = "NA"

$new.Range("B5").Value = "BD PROCUREMENT"
$new.Range("C5").Value = "Atleast by 20/09/2025 - for samples. 15/10/2025 - Large order"
$new.Range("D5").Value = "By today evening we will get the final date"

$new.Range("B6").Value = "NISE/ IIT LETTER"
$new.Range("C6").Value = "WILL BE COMMUNICATED SOON"

$new.Range("B7").Value = "IEC and PAN FILE with BD"
$new.Range("C7").Value = "WILL BE COMMUNICATED SOON"

# --- 3. Formatting -----------------------------------------------------------
# C3/C4 reuse the workbook's existing date-number-format style (same style
# already used for the date column on "Daily Updates"!B2).
$dateSource = $wb.Worksheets.Item("Daily Updates").Range("B2")
$dateSource.Copy()
$new.Range("C3").PasteSpecial(-4122)
$new.Range("C4").PasteSpecial(-4122)

# D3/C5 wrap their long text.
$new.Range("D3").WrapText = $true
$new.Range("C5").WrapText = $true

# Row heights for the two wrapped rows.
$new.Rows.Item(3).RowHeight = 58
$new.Rows.Item(5).RowHeight = 58

# Column widths.
$new.Columns.Item(2).ColumnWidth = 20.6328125
$new.Columns.Item(3).ColumnWidth = 27.90625
$new.Columns.Item(4).ColumnWidth = 36.36328125
$new.Columns.Item(5).ColumnWidth = 16

# --- 4. Selection / view state ------------------------------------------------
$new.Range("B2:E7").Select()
$excel.Application.CutCopyMode = $false

# Make the new sheet the active tab (matches activeTab="1" in the saved file).
$new.Activate()
